# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment Schedule"
#   sheet. This shifts the old N/O/P columns ("Late" / "heading" / "Outstanding")
#   one column to the right (N->O, O->P, P->Q).
# - Make "Repayment Schedule" the active/selected sheet (it was previously
#   "Transactions"), with a new selection.
# - "Transactions" loses its tabSelected flag and keeps its own selection.

$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsRepay.Columns("N").Insert() | Out-Null

$wsRepay.Activate() | Out-Null
$wsRepay.Range("L20").Select() | Out-Null
